# Apply template restoration edits to the "Open Action Items" sheet.
# Changes:
#  - Owner (col E) titles for rows tied to Chief Technology Officer / IT Managers /
#    DevOps Engineers / System Administrators are renamed to the Data/AI equivalents
#    for rows 8-11 and 14-17.
#  - Dependencies (col I) text for rows 8-17 updated from Cloud Infrastructure
#    Migration wording to AI/ML Implementation wording.
#  - Notes (col J) text for rows 8-17 updated from Information Technology wording
#    to Artificial Intelligence and Machine Learning wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Open Action Items")

# Owner (column E) updates
$ws.Range("E8").Value  = "Chief Data Officer"
$ws.Range("E9").Value  = "Data Scientists"
$ws.Range("E10").Value = "ML Engineers"
$ws.Range("E11").Value = "Business Analysts"
$ws.Range("E14").Value = "Chief Data Officer"
$ws.Range("E15").Value = "Data Scientists"
$ws.Range("E16").Value = "ML Engineers"
$ws.Range("E17").Value = "Business Analysts"

# Dependencies (column I) and Notes (column J) updates for rows 8-17
for ($row = 8; $row -le 17; $row++) {
    $ws.Cells.Item($row, 9).Value  = "Dependent on AI/ML Implementation milestone completion"
    $ws.Cells.Item($row, 10).Value = "Critical action for Artificial Intelligence and Machine Learning success"
}
